$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 156.656447
$ws.Range("H2").Value = 469.969341
$ws.Range("I2").Value = 0.0671576211124673
$ws.Range("J2").Value = 0.0671576211124673
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 3702.70008447089
$ws.Range("R2").Value = 33324.30076023801
$ws.Range("S2").Value = 0.004585417998302734
$ws.Range("T2").Value = 0.004585417998302734

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 156.656447
$ws.Range("H3").Value = 469.969341
$ws.Range("I3").Value = 0.0671576211124673
$ws.Range("J3").Value = 0.0671576211124673
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 28399.99526862559
$ws.Range("R3").Value = 255599.9574176303
$ws.Range("S3").Value = 0.03517050975925245
$ws.Range("T3").Value = 0.03517050975925244

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 156.656447
$ws.Range("H4").Value = 469.969341
$ws.Range("I4").Value = 0.0671576211124673
$ws.Range("J4").Value = 0.0671576211124673
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 17404.61684833888
$ws.Range("R4").Value = 156641.5516350499
$ws.Range("S4").Value = 0.02155385030633407
$ws.Range("T4").Value = 0.02155385030633407

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 156.656447
$ws.Range("H5").Value = 469.969341
$ws.Range("I5").Value = 0.0671576211124673
$ws.Range("J5").Value = 0.0671576211124673
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 4722.101443741248
$ws.Range("R5").Value = 42498.91299367123
$ws.Range("S5").Value = 0.005847843048578047
$ws.Range("T5").Value = 0.005847843048578047

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 42.300692
$ws.Range("H6").Value = 126.902076
$ws.Range("I6").Value = 0.01813403725498241
$ws.Range("J6").Value = 0.01813403725498241
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 999.8105972719852
$ws.Range("R6").Value = 8998.295375447868
$ws.Range("S6").Value = 0.001238163881231524
$ws.Range("T6").Value = 0.001238163881231524

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.300692
$ws.Range("H7").Value = 126.902076
$ws.Range("I7").Value = 0.01813403725498241
$ws.Range("J7").Value = 0.01813403725498241
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 7668.624404966803
$ws.Range("R7").Value = 69017.61964470122
$ws.Range("S7").Value = 0.009496812479151477
$ws.Range("T7").Value = 0.009496812479151476

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.300692
$ws.Range("H8").Value = 126.902076
$ws.Range("I8").Value = 0.01813403725498241
$ws.Range("J8").Value = 0.01813403725498241
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 4699.629991478062
$ws.Range("R8").Value = 42296.66992330256
$ws.Range("S8").Value = 0.005820014437212043
$ws.Range("T8").Value = 0.005820014437212043

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.300692
$ws.Range("H9").Value = 126.902076
$ws.Range("I9").Value = 0.01813403725498241
$ws.Range("J9").Value = 0.01813403725498241
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 1275.071422783218
$ws.Range("R9").Value = 11475.64280504897
$ws.Range("S9").Value = 0.001579046457387362
$ws.Range("T9").Value = 0.001579046457387362

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2110.189616
$ws.Range("H10").Value = 6330.568848
$ws.Range("I10").Value = 0.9046248489651427
$ws.Range("J10").Value = 0.9046248489651426
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 49876.01480207703
$ws.Range("R10").Value = 448884.1332186933
$ws.Range("S10").Value = 0.06176637878834273
$ws.Range("T10").Value = 0.06176637878834272

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2110.189616
$ws.Range("H11").Value = 6330.568848
$ws.Range("I11").Value = 0.9046248489651427
$ws.Range("J11").Value = 0.9046248489651426
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 382552.8808929445
$ws.Range("R11").Value = 3442975.9280365
$ws.Range("S11").Value = 0.4737528898724557
$ws.Range("T11").Value = 0.4737528898724555

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2110.189616
$ws.Range("H12").Value = 6330.568848
$ws.Range("I12").Value = 0.9046248489651427
$ws.Range("J12").Value = 0.9046248489651426
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 234443.2192045269
$ws.Range("R12").Value = 2109988.972840742
$ws.Range("S12").Value = 0.2903341162923514
$ws.Range("T12").Value = 0.2903341162923514

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2110.189616
$ws.Range("H13").Value = 6330.568848
$ws.Range("I13").Value = 0.9046248489651427
$ws.Range("J13").Value = 0.9046248489651426
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 63607.52859587956
$ws.Range("R13").Value = 572467.757362916
$ws.Range("S13").Value = 0.07877146401199295
$ws.Range("T13").Value = 0.07877146401199295

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 23.52144266666667
$ws.Range("H14").Value = 70.564328
$ws.Range("I14").Value = 0.01008349266740757
$ws.Range("J14").Value = 0.01008349266740757
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 555.9480596974338
$ws.Range("R14").Value = 5003.532537276904
$ws.Range("S14").Value = 0.0006884852083347659
$ws.Range("T14").Value = 0.0006884852083347657

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 23.52144266666667
$ws.Range("H15").Value = 70.564328
$ws.Range("I15").Value = 0.01008349266740757
$ws.Range("J15").Value = 0.01008349266740757
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 4264.164502879231
$ws.Range("R15").Value = 38377.48052591307
$ws.Range("S15").Value = 0.005280734656644531
$ws.Range("T15").Value = 0.005280734656644528

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 23.52144266666667
$ws.Range("H16").Value = 70.564328
$ws.Range("I16").Value = 0.01008349266740757
$ws.Range("J16").Value = 0.01008349266740757
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 2613.245130814844
$ws.Range("R16").Value = 23519.20617733359
$ws.Range("S16").Value = 0.003236238686214764
$ws.Range("T16").Value = 0.003236238686214763

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 23.52144266666667
$ws.Range("H17").Value = 70.564328
$ws.Range("I17").Value = 0.01008349266740757
$ws.Range("J17").Value = 0.01008349266740757
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 709.0077714780782
$ws.Range("R17").Value = 709.0077714780782
$ws.Range("S17").Value = 0.0008780341162135114
$ws.Range("T17").Value = 0.0008780341162135114
